# Add "Ancient Arena" timer logic: a new lookup column for the
# Ancient_Arena_930PM event, appended after the existing Demon_Gates_10PM
# column (column Q, the 17th column) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell with the event name (this also extends the used range /
# sheet dimension from A1:P3 to A1:Q3 automatically).
$ws.Range("Q1").Value = "Ancient_Arena_930PM"

# Give the new column a best-fit-style width so the header text isn't
# truncated, matching the sizing of the other header columns.
$ws.Columns.Item(17).ColumnWidth = 24.92

# Leave the cursor where it lands after entering the new header (one cell
# below what was just typed).
$ws.Range("Q2").Select() | Out-Null
